$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Fecha) date-serial updates for rows 2-133 ---
$dateUpdates = @(
    [pscustomobject]@{Row=2; Value=44239},
    [pscustomobject]@{Row=3; Value=44239},
    [pscustomobject]@{Row=4; Value=44285},
    [pscustomobject]@{Row=5; Value=44285},
    [pscustomobject]@{Row=6; Value=44267},
    [pscustomobject]@{Row=7; Value=44267},
    [pscustomobject]@{Row=8; Value=44313},
    [pscustomobject]@{Row=9; Value=44313},
    [pscustomobject]@{Row=10; Value=44266},
    [pscustomobject]@{Row=11; Value=44266},
    [pscustomobject]@{Row=12; Value=44383},
    [pscustomobject]@{Row=13; Value=44383},
    [pscustomobject]@{Row=14; Value=44442},
    [pscustomobject]@{Row=15; Value=44442},
    [pscustomobject]@{Row=16; Value=44237},
    [pscustomobject]@{Row=17; Value=44237},
    [pscustomobject]@{Row=18; Value=44398},
    [pscustomobject]@{Row=19; Value=44398},
    [pscustomobject]@{Row=20; Value=44420},
    [pscustomobject]@{Row=21; Value=44420},
    [pscustomobject]@{Row=22; Value=44462},
    [pscustomobject]@{Row=23; Value=44462},
    [pscustomobject]@{Row=24; Value=44194},
    [pscustomobject]@{Row=25; Value=44194},
    [pscustomobject]@{Row=26; Value=44341},
    [pscustomobject]@{Row=27; Value=44341},
    [pscustomobject]@{Row=28; Value=44460},
    [pscustomobject]@{Row=29; Value=44460},
    [pscustomobject]@{Row=30; Value=44371},
    [pscustomobject]@{Row=31; Value=44371},
    [pscustomobject]@{Row=32; Value=44299},
    [pscustomobject]@{Row=33; Value=44299},
    [pscustomobject]@{Row=34; Value=44217},
    [pscustomobject]@{Row=35; Value=44217},
    [pscustomobject]@{Row=36; Value=44365},
    [pscustomobject]@{Row=37; Value=44365},
    [pscustomobject]@{Row=38; Value=44405},
    [pscustomobject]@{Row=39; Value=44405},
    [pscustomobject]@{Row=40; Value=44376},
    [pscustomobject]@{Row=41; Value=44376},
    [pscustomobject]@{Row=42; Value=44222},
    [pscustomobject]@{Row=43; Value=44222},
    [pscustomobject]@{Row=44; Value=44257},
    [pscustomobject]@{Row=45; Value=44257},
    [pscustomobject]@{Row=46; Value=44327},
    [pscustomobject]@{Row=47; Value=44327},
    [pscustomobject]@{Row=48; Value=44278},
    [pscustomobject]@{Row=49; Value=44278},
    [pscustomobject]@{Row=50; Value=44245},
    [pscustomobject]@{Row=51; Value=44245},
    [pscustomobject]@{Row=52; Value=44336},
    [pscustomobject]@{Row=53; Value=44336},
    [pscustomobject]@{Row=54; Value=44271},
    [pscustomobject]@{Row=55; Value=44271},
    [pscustomobject]@{Row=56; Value=44308},
    [pscustomobject]@{Row=57; Value=44308},
    [pscustomobject]@{Row=58; Value=44330},
    [pscustomobject]@{Row=59; Value=44330},
    [pscustomobject]@{Row=60; Value=44322},
    [pscustomobject]@{Row=61; Value=44322},
    [pscustomobject]@{Row=62; Value=44224},
    [pscustomobject]@{Row=63; Value=44224},
    [pscustomobject]@{Row=64; Value=44264},
    [pscustomobject]@{Row=65; Value=44264},
    [pscustomobject]@{Row=66; Value=44209},
    [pscustomobject]@{Row=67; Value=44209},
    [pscustomobject]@{Row=68; Value=44447},
    [pscustomobject]@{Row=69; Value=44447},
    [pscustomobject]@{Row=70; Value=44316},
    [pscustomobject]@{Row=71; Value=44316},
    [pscustomobject]@{Row=72; Value=44169},
    [pscustomobject]@{Row=73; Value=44169},
    [pscustomobject]@{Row=74; Value=44467},
    [pscustomobject]@{Row=75; Value=44467},
    [pscustomobject]@{Row=76; Value=44427},
    [pscustomobject]@{Row=77; Value=44427},
    [pscustomobject]@{Row=78; Value=44203},
    [pscustomobject]@{Row=79; Value=44203},
    [pscustomobject]@{Row=80; Value=44469},
    [pscustomobject]@{Row=81; Value=44469},
    [pscustomobject]@{Row=82; Value=44168},
    [pscustomobject]@{Row=83; Value=44168},
    [pscustomobject]@{Row=84; Value=44434},
    [pscustomobject]@{Row=85; Value=44434},
    [pscustomobject]@{Row=86; Value=44292},
    [pscustomobject]@{Row=87; Value=44292},
    [pscustomobject]@{Row=88; Value=44161},
    [pscustomobject]@{Row=89; Value=44161},
    [pscustomobject]@{Row=90; Value=44280},
    [pscustomobject]@{Row=91; Value=44280},
    [pscustomobject]@{Row=92; Value=44274},
    [pscustomobject]@{Row=93; Value=44274},
    [pscustomobject]@{Row=94; Value=44344},
    [pscustomobject]@{Row=95; Value=44344},
    [pscustomobject]@{Row=96; Value=44358},
    [pscustomobject]@{Row=97; Value=44358},
    [pscustomobject]@{Row=98; Value=44391},
    [pscustomobject]@{Row=99; Value=44391},
    [pscustomobject]@{Row=100; Value=44231},
    [pscustomobject]@{Row=101; Value=44231},
    [pscustomobject]@{Row=102; Value=44320},
    [pscustomobject]@{Row=103; Value=44320},
    [pscustomobject]@{Row=104; Value=44475},
    [pscustomobject]@{Row=105; Value=44475},
    [pscustomobject]@{Row=106; Value=44252},
    [pscustomobject]@{Row=107; Value=44252},
    [pscustomobject]@{Row=108; Value=44204},
    [pscustomobject]@{Row=109; Value=44204},
    [pscustomobject]@{Row=110; Value=44362},
    [pscustomobject]@{Row=111; Value=44362},
    [pscustomobject]@{Row=112; Value=44350},
    [pscustomobject]@{Row=113; Value=44350},
    [pscustomobject]@{Row=114; Value=44453},
    [pscustomobject]@{Row=115; Value=44453},
    [pscustomobject]@{Row=116; Value=44435},
    [pscustomobject]@{Row=117; Value=44435},
    [pscustomobject]@{Row=118; Value=44433},
    [pscustomobject]@{Row=119; Value=44433},
    [pscustomobject]@{Row=120; Value=44159},
    [pscustomobject]@{Row=121; Value=44159},
    [pscustomobject]@{Row=122; Value=44166},
    [pscustomobject]@{Row=123; Value=44166},
    [pscustomobject]@{Row=124; Value=44334},
    [pscustomobject]@{Row=125; Value=44334},
    [pscustomobject]@{Row=126; Value=44386},
    [pscustomobject]@{Row=127; Value=44386},
    [pscustomobject]@{Row=128; Value=44306},
    [pscustomobject]@{Row=129; Value=44306},
    [pscustomobject]@{Row=130; Value=44425},
    [pscustomobject]@{Row=131; Value=44425},
    [pscustomobject]@{Row=132; Value=44187},
    [pscustomobject]@{Row=133; Value=44187}
)

foreach ($u in $dateUpdates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.Value
}

# --- Column J (Volumen) updates ---
$ws.Cells.Item(22, 10).Value = 200
$ws.Cells.Item(23, 10).Value = 100
$ws.Cells.Item(116, 10).Value = 400
$ws.Cells.Item(117, 10).Value = 200

# --- Column O (Origen) updates ---
$ws.Cells.Item(38, 15).Value = "Región de Ñuble"
$ws.Cells.Item(39, 15).Value = "Región de Ñuble"
$ws.Cells.Item(70, 15).Value = "Región Metropolitana"
$ws.Cells.Item(71, 15).Value = "Región Metropolitana"
$ws.Cells.Item(102, 15).Value = "Región Metropolitana"
$ws.Cells.Item(103, 15).Value = "Región Metropolitana"
$ws.Cells.Item(112, 15).Value = "Región de Ñuble"
$ws.Cells.Item(113, 15).Value = "Región de Ñuble"

# --- New rows 134 and 135 ---
$ws.Cells.Item(134, 1).Value = 11
$ws.Cells.Item(134, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(134, 3).Value = "Bíobío"
$ws.Cells.Item(134, 4).Value = 44250
$ws.Cells.Item(134, 4).NumberFormat = $ws.Cells.Item(133, 4).NumberFormat
$ws.Cells.Item(134, 5).Value = 8
$ws.Cells.Item(134, 6).Value = 100112040
$ws.Cells.Item(134, 7).Value = "Cilantro"
$ws.Cells.Item(134, 8).Value = "Sin especificar"
$ws.Cells.Item(134, 9).Value = "Primera"
$ws.Cells.Item(134, 10).Value = 200
$ws.Cells.Item(134, 11).Value = 600
$ws.Cells.Item(134, 12).Value = 700
$ws.Cells.Item(134, 13).Value = 650
$ws.Cells.Item(134, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(134, 15).Value = "Región de Ñuble"
$ws.Cells.Item(134, 16).Value = 650
$ws.Cells.Item(134, 17).Value = 1
$ws.Cells.Item(134, 18).Value = "Hortaliza"

$ws.Cells.Item(135, 1).Value = 11
$ws.Cells.Item(135, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(135, 3).Value = "Bíobío"
$ws.Cells.Item(135, 4).Value = 44250
$ws.Cells.Item(135, 4).NumberFormat = $ws.Cells.Item(133, 4).NumberFormat
$ws.Cells.Item(135, 5).Value = 8
$ws.Cells.Item(135, 6).Value = 100112040
$ws.Cells.Item(135, 7).Value = "Cilantro"
$ws.Cells.Item(135, 8).Value = "Sin especificar"
$ws.Cells.Item(135, 9).Value = "Segunda"
$ws.Cells.Item(135, 10).Value = 100
$ws.Cells.Item(135, 11).Value = 500
$ws.Cells.Item(135, 12).Value = 500
$ws.Cells.Item(135, 13).Value = 500
$ws.Cells.Item(135, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(135, 15).Value = "Región de Ñuble"
$ws.Cells.Item(135, 16).Value = 500
$ws.Cells.Item(135, 17).Value = 1
$ws.Cells.Item(135, 18).Value = "Hortaliza"
